# Updated symbol list on Sat Feb 11 03:58:43 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# the crypto rows on the sheet to the latest scraped figures. All values in
# these columns are stored as plain text (not numbers/percentages), so we
# force the cell format to Text before writing and then clear the format
# again so no stray number-format style is left behind on the cell.

function Set-CellText($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "308.05"
Set-CellText $ws "E2" "0.16%"
Set-CellText $ws "D3" "41.31"
Set-CellText $ws "E3" "3.08%"
Set-CellText $ws "D4" "5.130"
Set-CellText $ws "E4" "2.34%"
Set-CellText $ws "D5" "0.07613"
Set-CellText $ws "E5" "-0.76%"
Set-CellText $ws "D6" "1.623"
Set-CellText $ws "E6" "0.18%"
Set-CellText $ws "D8" "0.9037"
Set-CellText $ws "E8" "1.92%"
Set-CellText $ws "D9" "0.1094"
Set-CellText $ws "E9" "8.88%"
Set-CellText $ws "D10" "0.1766"
Set-CellText $ws "E10" "1.87%"
Set-CellText $ws "D11" "0.09159"
Set-CellText $ws "E11" "2.37%"
Set-CellText $ws "D12" "0.04273"
Set-CellText $ws "E12" "-2.57%"
Set-CellText $ws "E13" "-0.50%"
Set-CellText $ws "D14" "0.001251"
Set-CellText $ws "E14" "-1.14%"
Set-CellText $ws "D15" "0.005824"
Set-CellText $ws "E15" "0.59%"
Set-CellText $ws "D16" "3.360"
Set-CellText $ws "E16" "0.10%"
Set-CellText $ws "D17" "4.251"
Set-CellText $ws "E17" "0.62%"
Set-CellText $ws "D19" "6.569"
Set-CellText $ws "E19" "-6.54%"
Set-CellText $ws "D20" "0.1360"
Set-CellText $ws "E20" "1.32%"
Set-CellText $ws "D22" "0.04161"
Set-CellText $ws "E22" "-1.85%"
Set-CellText $ws "D23" "0.001224"
Set-CellText $ws "E23" "1.94%"
Set-CellText $ws "D24" "0.004084"
Set-CellText $ws "E24" "0.38%"
Set-CellText $ws "D25" "0.0001300"
Set-CellText $ws "E25" "6.39%"
Set-CellText $ws "D38" "0.02414"
Set-CellText $ws "E38" "2.76%"
Set-CellText $ws "E39" "0.89%"
Set-CellText $ws "D40" "0.007754"
Set-CellText $ws "E40" "-2.61%"
Set-CellText $ws "D41" "0.1298"
Set-CellText $ws "E41" "-1.92%"
Set-CellText $ws "D42" "0.006956"
Set-CellText $ws "E42" "5.85%"
Set-CellText $ws "D43" "0.001920"
Set-CellText $ws "E43" "-4.10%"
Set-CellText $ws "D44" "0.008076"
Set-CellText $ws "E44" "6.00%"
Set-CellText $ws "D45" "0.3056"
Set-CellText $ws "E45" "0.18%"
Set-CellText $ws "D46" "0.00006743"
Set-CellText $ws "E46" "2.38%"
Set-CellText $ws "E47" "-0.16%"
Set-CellText $ws "D48" "0.009298"
Set-CellText $ws "E48" "178.57%"
Set-CellText $ws "D49" "0.004203"
Set-CellText $ws "E49" "-15.99%"
Set-CellText $ws "E50" "-0.16%"
Set-CellText $ws "E51" "-0.16%"
